# C5-PowerPoint.pptx edit: swap the deck's colour theme for the default
# "Office Theme" palette (was the "Integral" theme) and pick a different
# built-in table style for the sources-of-finance table.

$p = $ppt.ActivePresentation

# --- 1. Table style on the "SOURCES OF FINANCE" table (slide 6) -----------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{6BBA0C2D-ACB5-47AA-9CF3-216FC37C46A3}", $true)
    }
}

# --- 2. Theme colours: Integral -> Office Theme ----------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (RGB() == r + g*256 + b*65536)
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
